# Update port arrivals Excel
# - Refresh "Last Updated" timestamps on Sydney and Melbourne sheets
# - Remove the now-redundant "Last Updated" summary sheet

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2026-02-18 00:28"
$newTimestamp = "2026-02-18 00:40"

$sydney = $wb.Worksheets.Item("Sydney")
$lastRowSydney = $sydney.Cells.Item($sydney.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowSydney; $r++) {
    $cell = $sydney.Cells.Item($r, 8)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

$melbourne = $wb.Worksheets.Item("Melbourne")
$lastRowMelbourne = $melbourne.Cells.Item($melbourne.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowMelbourne; $r++) {
    $cell = $melbourne.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

$lastUpdated = $wb.Worksheets.Item("Last Updated")
$lastUpdated.Delete()
